# Update football match bases for "Estonia Meistriliiga".
#
# The source data rows got re-ordered: for each group of rows listed below,
# the record data (columns B..AC - id, teams, scores, odds, etc.) is rotated
# among the rows of the group while the row's own sequence number in column A
# stays put. For a 2-row group this is simply a swap; for the 4-row group at
# the end it is a 4-cycle rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-RowData($Worksheet, $Rows) {
    $count = $Rows.Length

    # Snapshot the current B:AC contents of every row in the group first,
    # so overwriting one row never affects the data we still need to read
    # from another row.
    $ranges = @()
    $values = @()
    for ($i = 0; $i -lt $count; $i++) {
        $r = $Rows[$i]
        $rng = $Worksheet.Range("B$r`:AC$r")
        $ranges += $rng
        $values += $rng.Value()
    }

    # Write each row's new data from the row that precedes it in the group
    # (wrapping around), i.e. rotate the snapshotted data by one position.
    for ($i = 0; $i -lt $count; $i++) {
        $srcIndex = ($i - 1 + $count) % $count
        $ranges[$i].Value = $values[$srcIndex]
    }
}

Rotate-RowData $ws @(32, 33)
Rotate-RowData $ws @(40, 41)
Rotate-RowData $ws @(77, 78)
Rotate-RowData $ws @(161, 162)
Rotate-RowData $ws @(168, 169)
Rotate-RowData $ws @(177, 178, 179, 180)
